# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet, shifting the "Late" / "Outstanding" columns one place to the
# right, and update the active selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new entire column before column N (shifts N:P -> O:Q)
$ws.Columns("N").Insert()

# New column is a bit wider than the default bestFit width used by its
# neighbours. Excel stores width in "characters" while padding on an
# internal measure, so request a value that round-trips to exactly 10.
$ws.Columns("N").ColumnWidth = 9.166666666666666

# Update the selected cell/range recorded for this sheet
$ws.Range("T9").Select()
